$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (interested count) column F
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1178
$ws.Range("F4").Value = 1612
$ws.Range("F5").Value = 185
$ws.Range("F6").Value = 185
$ws.Range("F7").Value = 40
$ws.Range("F8").Value = 1565
$ws.Range("F9").Value = 3174
$ws.Range("F10").Value = 715
$ws.Range("F11").Value = 1899
$ws.Range("F12").Value = 1854
$ws.Range("F13").Value = 925
$ws.Range("F14").Value = 317
$ws.Range("F16").Value = 1539
$ws.Range("F17").Value = 312
$ws.Range("F20").Value = 1335
$ws.Range("F21").Value = 445
$ws.Range("F22").Value = 548
$ws.Range("F23").Value = 235
$ws.Range("F24").Value = 7930
$ws.Range("F25").Value = 9268
$ws.Range("F26").Value = 797
$ws.Range("F27").Value = 610
$ws.Range("F28").Value = 1757
$ws.Range("F30").Value = 298

# Sheet "演出" (performances) - update column F
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 29

# Sheet "全部类型" (all types) - update column F
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1178
$ws.Range("F5").Value = 1612
$ws.Range("F6").Value = 185
$ws.Range("F7").Value = 185
$ws.Range("F9").Value = 40
$ws.Range("F10").Value = 1565
$ws.Range("F11").Value = 3174
$ws.Range("F12").Value = 715
$ws.Range("F13").Value = 1899
$ws.Range("F14").Value = 1854
$ws.Range("F15").Value = 925
$ws.Range("F16").Value = 317
$ws.Range("F18").Value = 1539
$ws.Range("F19").Value = 312
$ws.Range("F24").Value = 1335
$ws.Range("F25").Value = 445
$ws.Range("F26").Value = 548
$ws.Range("F27").Value = 235
$ws.Range("F28").Value = 7930
$ws.Range("F29").Value = 9268
$ws.Range("F30").Value = 797
$ws.Range("F31").Value = 610
$ws.Range("F32").Value = 1757
$ws.Range("F34").Value = 29
$ws.Range("F36").Value = 298
